$wb = $excel.ActiveWorkbook

# sheet1 (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 37609
$ws.Range("F6").Value = 476
$ws.Range("F7").Value = 363
$ws.Range("F10").Value = 93
$ws.Range("F11").Value = 711
$ws.Range("F13").Value = 36
$ws.Range("F15").Value = 16
$ws.Range("F18").Value = 469
$ws.Range("F20").Value = 1165
$ws.Range("F22").Value = 826
$ws.Range("F23").Value = 2521
$ws.Range("F24").Value = 1005
$ws.Range("F25").Value = 562
$ws.Range("F26").Value = 107
$ws.Range("F27").Value = 1158
$ws.Range("F29").Value = 769
$ws.Range("F30").Value = 57
$ws.Range("F31").Value = 1156

# sheet2 (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 394
$ws.Range("F5").Value = 3
$ws.Range("F9").Value = 143

# sheet3 (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 629

# sheet4 (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 629
$ws.Range("F3").Value = 37609
$ws.Range("F7").Value = 476
$ws.Range("F9").Value = 363
$ws.Range("F11").Value = 394
$ws.Range("F13").Value = 3
$ws.Range("F16").Value = 93
$ws.Range("F17").Value = 711
$ws.Range("F19").Value = 36
$ws.Range("F23").Value = 143
$ws.Range("F25").Value = 16
$ws.Range("F29").Value = 469
$ws.Range("F31").Value = 1165
$ws.Range("F33").Value = 826
$ws.Range("F34").Value = 2521
$ws.Range("F35").Value = 1005
$ws.Range("F36").Value = 562
$ws.Range("F37").Value = 107
$ws.Range("F38").Value = 1158
$ws.Range("F41").Value = 769
$ws.Range("F42").Value = 57
$ws.Range("F43").Value = 1156

